# Updated cryptos list on Sat Aug 12 15:51:13 UTC 2023 with GitHub Actions
# Applies latest scraped price / volume(1h) data to the cryptos sheet,
# including two rows (14/15) whose Polygon / ShibaInu entries swapped rank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.464.57"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "1.852.35"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6306"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07656"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2939"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.64"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07758"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.81%  "
$ws.Range("D12").Value = "1.862.16"
$ws.Range("E12").Value = "  +0.94%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.032"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6806"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001066"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.66%  "
$ws.Range("D17").Value = "2.107.04"
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.170"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "29.478.14"
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "229.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.26%  "
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.448"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.002"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("E26").Value = "  -0.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.403"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.326"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.469"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05684"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.133"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.045"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.850"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.72%  "
$ws.Range("E35").Value = "  +1.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7023"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.95%  "
$ws.Range("E37").Value = "  -0.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.785"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.43%  "
$ws.Range("E39").Value = "  -0.61%  "
$ws.Range("D40").Value = "1.218.11"
$ws.Range("E40").Value = "  -2.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.552"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9089"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.002"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("D44").Value = "2.015.71"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.65"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "66.46"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.33%  "
$ws.Range("E47").Value = "  -1.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.106"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4020"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.006"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.683"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.03%  "
